# v2.0.05061200 long press send pulse
# Adds a new column F of pulse-length values alongside the existing A:E data,
# and moves the active selection to the newly added cell F20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value  = 0
$ws.Range("F4").Value  = 20
$ws.Range("F6").Value  = 24
$ws.Range("F16").Value = 2
$ws.Range("F20").Value = 0

[void]$ws.Range("F20").Select()
